# Add presenter name + restyle the "Last updated" subtitle placeholder on slide 1.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)   # "副标题 4" - the subTitle placeholder holding "Last updated: 9/9/2019"

# Give the placeholder an explicit position/size (it previously inherited the
# layout's <p:spPr/> with no xfrm).
$sh.Left = 120
$sh.Top = 423.2704724409449
$sh.Width = 720
$sh.Height = 29.624961853027344

# Stop PowerPoint from auto-shrinking the text to fit (normAutofit -> noAutofit).
$sh.TextFrame.AutoSize = 0

$tr = $sh.TextFrame.TextRange

# Insert a new first paragraph with the author's name, ahead of the existing
# "Last updated: 9/9/2019" paragraph (which keeps its own formatting/endParaRPr).
[void]$tr.InsertBefore("Sparks Lu`r")

# Both paragraphs now render at an explicit 20pt (sz="2000") instead of relying
# on the placeholder's default/autofit-scaled size.
$tr.Font.Size = 20
